# "more fixes to the slide masters"
#
# Removes the stray/duplicate "Slide Number Placeholder" shapes that were
# left behind on the slide master and on every slide layout.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# --- Slide layouts: each one has exactly one slide-number placeholder ---
# (idx="12") that needs to go. The shape's generated name varies per
# layout ("Slide Number Placeholder 3/4/5/6/8") but is unique within the
# layout, so look it up by name.
$layoutPlaceholderNames = @(
    "Slide Number Placeholder 5",  # 1  Title Slide
    "Slide Number Placeholder 5",  # 2  Title and Content
    "Slide Number Placeholder 5",  # 3  Section Header
    "Slide Number Placeholder 6",  # 4  Two Content
    "Slide Number Placeholder 8",  # 5  Comparison
    "Slide Number Placeholder 4",  # 6  Title Only
    "Slide Number Placeholder 3",  # 7  Blank
    "Slide Number Placeholder 6",  # 8  Content with Caption
    "Slide Number Placeholder 6",  # 9  Picture with Caption
    "Slide Number Placeholder 5",  # 10 Title and Vertical Text
    "Slide Number Placeholder 5"   # 11 Vertical Title and Text
)

for ($li = 1; $li -le $layoutPlaceholderNames.Count; $li++) {
    $lo = $m.CustomLayouts.Item($li)
    $name = $layoutPlaceholderNames[$li - 1]
    $sh = $lo.Shapes.Item($name)
    $sh.Delete()
}

# --- Slide master itself ---
# The master has two shapes both (confusingly) named
# "Slide Number Placeholder 5": the real placeholder (idx="4", the one to
# remove) and a userDrawn text box duplicate (idx is absent) that must stay.
# Match on the placeholder's index within the master's Shapes collection
# rather than by name to avoid ambiguity.
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Name -eq "Slide Number Placeholder 5" -and $sh.PlaceholderFormat.Type -eq 13) {
        $sh.Delete()
        break
    }
}
